$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.594256754476419
$ws.Range("C2").Value = 0.2068058676653948
$ws.Range("E2").Value = 0.08464720699963557
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.9996598946632673
$ws.Range("H2").Value = 0.9948011047668928
$ws.Range("L2").Value = 0.2284930407219719
$ws.Range("B3").Value = 1.467279717020176
$ws.Range("C3").Value = 0.1895804479934782
$ws.Range("E3").Value = 0.08525993122890085
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.9972475033090973
$ws.Range("H3").Value = 1.001224761628606
$ws.Range("L3").Value = 0.2188501910666218
$ws.Range("B4").Value = 1.389840441853039
$ws.Range("C4").Value = 0.1789269812140901
$ws.Range("E4").Value = 0.08567575437267116
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.996835350407963
$ws.Range("H4").Value = 1.005927678669892
$ws.Range("L4").Value = 0.2130565449390787
$ws.Range("B5").Value = 1.358415599380976
$ws.Range("C5").Value = 0.1745663253504404
$ws.Range("E5").Value = 0.08585516233785206
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.9969343464569249
$ws.Range("H5").Value = 1.008034223351146
$ws.Range("L5").Value = 0.2107274742437113
$ws.Range("B6").Value = 1.353205526329418
$ws.Range("C6").Value = 0.1738410796224059
$ws.Range("E6").Value = 0.08588555416254273
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.9969668530795701
$ws.Range("H6").Value = 1.008395471770157
$ws.Range("L6").Value = 0.2103426576496048
$ws.Range("B7").Value = 1.389416099450102
$ws.Range("C7").Value = 0.1788682498439869
$ws.Range("E7").Value = 0.08567813362120269
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.9968356071144058
$ws.Range("H7").Value = 1.005955319669567
$ws.Range("L7").Value = 0.2130250052612297
$ws.Range("B8").Value = 1.550366169748315
$ws.Range("C8").Value = 0.2008825623651944
$ws.Range("E8").Value = 0.08485025177770389
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.9986050864754503
$ws.Range("H8").Value = 0.9968580591419993
$ws.Range("L8").Value = 0.2251417572213086
$ws.Range("B9").Value = 1.870164278908419
$ws.Range("C9").Value = 0.2434404998932109
$ws.Range("E9").Value = 0.08354121184072483
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 1.010643705149562
$ws.Range("H9").Value = 0.9850723688647633
$ws.Range("L9").Value = 0.2499161303019974
$ws.Range("B10").Value = 2.107705893413765
$ws.Range("C10").Value = 0.2743364079192645
$ws.Range("E10").Value = 0.08277140586577048
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 1.024833746467522
$ws.Range("H10").Value = 0.9801514494415926
$ws.Range("L10").Value = 0.268744908329424
$ws.Range("B11").Value = 2.216342346670274
$ws.Range("C11").Value = 0.2883117726792079
$ws.Range("E11").Value = 0.08246294090225703
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 1.032476204802606
$ws.Range("H11").Value = 0.9787347745953525
$ws.Range("L11").Value = 0.2774489424064086
$ws.Range("B12").Value = 2.257563597434967
$ws.Range("C12").Value = 0.2935924650618631
$ws.Range("E12").Value = 0.0823521381558443
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 1.035543042637613
$ws.Range("H12").Value = 0.9783173326215433
$ws.Range("L12").Value = 0.2807650170843345
$ws.Range("B13").Value = 2.248682173484838
$ws.Range("C13").Value = 0.292455684814712
$ws.Range("E13").Value = 0.0823757342770719
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 1.03487482571046
$ws.Range("H13").Value = 0.9784019293028337
$ws.Range("L13").Value = 0.2800499479447893
$ws.Range("B14").Value = 2.219731980915753
$ws.Range("C14").Value = 0.2887464487525335
$ws.Range("E14").Value = 0.0824537046823508
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 1.032725040268645
$ws.Range("H14").Value = 0.9786980412038417
$ws.Range("L14").Value = 0.2777213555490903
$ws.Range("B15").Value = 2.202009957599728
$ws.Range("C15").Value = 0.2864729383244082
$ws.Range("E15").Value = 0.08250224618958057
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 1.031430800351501
$ws.Range("H15").Value = 0.9788949433939251
$ws.Range("L15").Value = 0.2762976389712151
$ws.Range("B16").Value = 2.100617860953321
$ws.Range("C16").Value = 0.2734214826780601
$ws.Range("E16").Value = 0.08279240491027373
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 1.024358362650844
$ws.Range("H16").Value = 0.9802606422009603
$ws.Range("L16").Value = 0.2681788790066406
$ws.Range("B17").Value = 2.038564960898441
$ws.Range("C17").Value = 0.2653944678737048
$ws.Range("E17").Value = 0.08298109885406291
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 1.020325225425523
$ws.Range("H17").Value = 0.9813095508057188
$ws.Range("L17").Value = 0.2632338822564435
$ws.Range("B18").Value = 2.002928140200027
$ws.Range("C18").Value = 0.2607700770974191
$ws.Range("E18").Value = 0.08309355745585378
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 1.018117105448013
$ws.Range("H18").Value = 0.981990163888895
$ws.Range("L18").Value = 0.2604027026912661
$ws.Range("B19").Value = 1.990871465402563
$ws.Range("C19").Value = 0.2592030579758671
$ws.Range("E19").Value = 0.08313230824950857
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 1.01738858430457
$ws.Range("H19").Value = 0.9822338601534568
$ws.Range("L19").Value = 0.2594463526972532
$ws.Range("B20").Value = 2.045164972009104
$ws.Range("C20").Value = 0.2662497303061855
$ws.Range("E20").Value = 0.08296060560301299
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 1.020742991157334
$ws.Range("H20").Value = 0.9811898855237757
$ws.Range("L20").Value = 0.263758934100025
$ws.Range("B21").Value = 2.228233102434558
$ws.Range("C21").Value = 0.2898362532223473
$ws.Range("E21").Value = 0.08243063982486554
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 1.033351777985985
$ws.Range("H21").Value = 0.9786078290428577
$ws.Range("L21").Value = 0.2784047747531133
$ws.Range("B22").Value = 2.348362399163022
$ws.Range("C22").Value = 0.3051844982113039
$ws.Range("E22").Value = 0.08211928858399631
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 1.042600551547253
$ws.Range("H22").Value = 0.9776144489894989
$ws.Range("L22").Value = 0.2880935684565031
$ws.Range("B23").Value = 2.284202887031825
$ws.Range("C23").Value = 0.2969989936916306
$ws.Range("E23").Value = 0.08228225684312562
$ws.Range("F23").Value = 0.7472568307915566
$ws.Range("G23").Value = 1.037571371634414
$ws.Range("H23").Value = 0.9780808427348404
$ws.Range("L23").Value = 0.2829117485195951
$ws.Range("B24").Value = 2.042180987939673
$ws.Range("C24").Value = 0.2658630959149093
$ws.Range("E24").Value = 0.0829698582180125
$ws.Range("F24").Value = 0.6416283278902313
$ws.Range("G24").Value = 1.020553774969443
$ws.Range("H24").Value = 0.9812437445979754
$ws.Range("L24").Value = 0.2635215215250781
$ws.Range("B25").Value = 1.783199384725606
$ws.Range("C25").Value = 0.2319928200656705
$ws.Range("E25").Value = 0.08386165315665473
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 1.006457056464754
$ws.Range("H25").Value = 0.987607705835714
$ws.Range("L25").Value = 0.2431045683398594
